$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.933
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("B23").Value = 7.886
$ws.Range("B25").Value = 5.392999999999999
$ws.Range("B53").Value = 6.412000000000001
$ws.Range("B57").Value = 5.034000000000001
$ws.Range("B59").Value = 4.682
$ws.Range("B69").Value = 5.59
$ws.Range("B79").Value = 5.577
$ws.Range("B83").Value = 5.824
$ws.Range("B93").Value = 4.973
